$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1333.3334
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1333.3334
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1683.3334

# Sheet ALC, row 64 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3296.7646
$ws.Range("I64").Value = 3202
$ws.Range("J64").Value = 3470.5
$ws.Range("K64").Value = 3202
$ws.Range("L64").Value = 3470.5
$ws.Range("M64").Value = -2954
$ws.Range("N64").Value = -3966.5

# Sheet ALC, row 67 (Leve Item ID 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3296.7646
$ws.Range("I67").Value = 3202
$ws.Range("J67").Value = 3470.5
$ws.Range("K67").Value = 3202
$ws.Range("L67").Value = 3470.5
$ws.Range("M67").Value = -2344
$ws.Range("N67").Value = -5186.5

# Sheet ALC, row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2472717.5
$ws.Range("I76").Value = 2852477.5
$ws.Range("K76").Value = 2852477.5
$ws.Range("M76").Value = -2852162.5

# Sheet ALC, row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2472717.5
$ws.Range("I79").Value = 2852477.5
$ws.Range("K79").Value = 2852477.5
$ws.Range("M79").Value = -2851385.5

# Sheet ALC, row 116 (Leve Item ID 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8994.727999999999
$ws.Range("I116").Value = 13100
$ws.Range("J116").Value = 4889.4546
$ws.Range("K116").Value = 13100
$ws.Range("L116").Value = 4889.4546
$ws.Range("M116").Value = -9658
$ws.Range("N116").Value = -11773.4546

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1234.5
$ws.Range("I137").Value = 1081.091
$ws.Range("K137").Value = 3243.273
$ws.Range("M137").Value = -693.2729999999997

# Sheet ALC, row 139 (Leve Item ID 42306)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 69772.5
$ws.Range("J139").Value = 69772.5
$ws.Range("L139").Value = 69772.5
$ws.Range("N139").Value = -80052.5

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3588.4194
$ws.Range("I61").Value = 3591.3667
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 3591.3667
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -3379.3667
$ws.Range("N61").Value = -3924

# Sheet ARM, row 63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3655.5557
$ws.Range("I63").Value = 2560
$ws.Range("J63").Value = 5025
$ws.Range("K63").Value = 2560
$ws.Range("L63").Value = 5025
$ws.Range("M63").Value = -1874
$ws.Range("N63").Value = -6397

# Sheet ARM, row 66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3655.5557
$ws.Range("I66").Value = 2560
$ws.Range("J66").Value = 5025
$ws.Range("K66").Value = 12800
$ws.Range("L66").Value = 25125
$ws.Range("M66").Value = -9368
$ws.Range("N66").Value = -31989

# Sheet ARM, row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 596.25
$ws.Range("I74").Value = 596.25
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 596.25
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 277.75
$ws.Range("N74").ClearContents()

# Sheet ARM, row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 596.25
$ws.Range("I77").Value = 596.25
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2981.25
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1386.75
$ws.Range("N77").ClearContents()

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3588.4194
$ws.Range("I136").Value = 3591.3667
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 10774.1001
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -8224.1001
$ws.Range("N136").Value = -15600

# Sheet ARM, row 139 (Leve Item ID 42321)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 64633.332
$ws.Range("J139").Value = 64633.332
$ws.Range("L139").Value = 64633.332
$ws.Range("N139").Value = -74913.33199999999

# Sheet ARM, row 141 (Leve Item ID 42483)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 62307.69
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 62307.69
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 62307.69
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -72667.69

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2163.805
$ws.Range("I86").Value = 1776.762
$ws.Range("J86").Value = 2570.2
$ws.Range("K86").Value = 1776.762
$ws.Range("L86").Value = 2570.2
$ws.Range("M86").Value = -653.7619999999999
$ws.Range("N86").Value = -4816.2

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2163.805
$ws.Range("I89").Value = 1776.762
$ws.Range("J89").Value = 2570.2
$ws.Range("K89").Value = 8883.809999999999
$ws.Range("L89").Value = 12851
$ws.Range("M89").Value = -3267.809999999999
$ws.Range("N89").Value = -24083

# Sheet BSM, row 138 (Leve Item ID 42308)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 59300
$ws.Range("J138").Value = 59300
$ws.Range("L138").Value = 59300
$ws.Range("N138").Value = -69580

# Sheet BSM, row 140 (Leve Item ID 42471)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 89663.336
$ws.Range("J140").Value = 89663.336
$ws.Range("L140").Value = 89663.336
$ws.Range("N140").Value = -100023.336

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5473.3115
$ws.Range("I31").Value = 4456
$ws.Range("J31").Value = 6133.189
$ws.Range("K31").Value = 4456
$ws.Range("L31").Value = 6133.189
$ws.Range("M31").Value = -4161
$ws.Range("N31").Value = -6723.189

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5473.3115
$ws.Range("I34").Value = 4456
$ws.Range("J34").Value = 6133.189
$ws.Range("K34").Value = 4456
$ws.Range("L34").Value = 6133.189
$ws.Range("M34").Value = -4254
$ws.Range("N34").Value = -6537.189

# Sheet CRP, row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1125.8438
$ws.Range("I58").Value = 706.5625
$ws.Range("J58").Value = 1545.125
$ws.Range("K58").Value = 706.5625
$ws.Range("L58").Value = 1545.125
$ws.Range("M58").Value = -503.5625
$ws.Range("N58").Value = -1951.125

# Sheet CRP, row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1125.8438
$ws.Range("I136").Value = 706.5625
$ws.Range("J136").Value = 1545.125
$ws.Range("K136").Value = 2119.6875
$ws.Range("L136").Value = 4635.375
$ws.Range("M136").Value = 430.3125
$ws.Range("N136").Value = -9735.375

# Sheet GSM, row 80 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2646
$ws.Range("I80").Value = 2493.6
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 2493.6
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -1495.6
$ws.Range("N80").Value = -4896

# Sheet GSM, row 83 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2646
$ws.Range("I83").Value = 2493.6
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 12468
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -7476
$ws.Range("N83").Value = -24484

# Sheet GSM, row 138 (Leve Item ID 42325)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 63306.25
$ws.Range("J138").Value = 63306.25
$ws.Range("L138").Value = 63306.25
$ws.Range("N138").Value = -73586.25

# Sheet GSM, row 140 (Leve Item ID 42458)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 99873
$ws.Range("J140").Value = 99873
$ws.Range("L140").Value = 99873
$ws.Range("N140").Value = -110233

# Sheet GSM, row 141 (Leve Item ID 42504)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 66551.39999999999
$ws.Range("J141").Value = 66551.39999999999
$ws.Range("L141").Value = 66551.39999999999
$ws.Range("N141").Value = -76911.39999999999

# Sheet LTW, row 82 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2297.1765
$ws.Range("I82").Value = 2705.7778
$ws.Range("J82").Value = 1837.5
$ws.Range("K82").Value = 2705.7778
$ws.Range("L82").Value = 1837.5
$ws.Range("M82").Value = -2344.7778
$ws.Range("N82").Value = -2559.5

# Sheet LTW, row 85 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2297.1765
$ws.Range("I85").Value = 2705.7778
$ws.Range("J85").Value = 1837.5
$ws.Range("K85").Value = 2705.7778
$ws.Range("L85").Value = 1837.5
$ws.Range("M85").Value = -1457.7778
$ws.Range("N85").Value = -4333.5

# Sheet LTW, row 138 (Leve Item ID 42334)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 50895.7
$ws.Range("J138").Value = 50895.7
$ws.Range("L138").Value = 50895.7
$ws.Range("N138").Value = -61175.7

# Sheet LTW, row 139 (Leve Item ID 43310)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 64900
$ws.Range("J139").Value = 64900
$ws.Range("L139").Value = 64900
$ws.Range("N139").Value = -75180

# Sheet WVR, row 138 (Leve Item ID 42347)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 58300
$ws.Range("J138").Value = 58300
$ws.Range("L138").Value = 58300
$ws.Range("N138").Value = -68580

# Sheet WVR, row 139 (Leve Item ID 43312)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 53942.855
$ws.Range("J139").Value = 53942.855
$ws.Range("L139").Value = 53942.855
$ws.Range("N139").Value = -64222.855
